# Swap the values of column C (codeforiati:group-code) and column D
# (codeforiati:group-name) for every row of data, including the header.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
if ($lastRow -lt 1) { $lastRow = $ws.UsedRange.Rows.Count }

for ($r = 1; $r -le $lastRow; $r++) {
    $cellC = $ws.Cells.Item($r, 3)
    $cellD = $ws.Cells.Item($r, 4)
    $valC = $cellC.Value()
    $valD = $cellD.Value()
    $cellC.Value = $valD
    $cellD.Value = $valC
}
